# Sprint Backlogs.xlsx — add the "Add uniform formatting to documents" backlog
# item (Sprint 1 table, row 11) and give its Task cell (C11) a grey "card"
# border, matching the commit "Moved product backlog into sperate excel file
# and added formatting to Project Plan document".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New backlog row -------------------------------------------------
$ws.Range("C11").Value = "Add uniform formatting to documents"
$ws.Range("D11").Value = "Morgan"
$ws.Range("E11").Value = "N/A"
$ws.Range("F11").Value = "1 hour"

# --- Formatting: medium grey border around the Task cell (left/right/bottom) -
$taskCell = $ws.Range("C11")
$taskCell.Borders.Weight = -4138          # xlMedium, applied to all 4 sides first
$taskCell.Borders.Color = 10066329        # RGB(153,153,153) = FF999999
$taskCell.Borders.Item(8).LineStyle = -4142  # xlEdgeTop -> xlLineStyleNone (no top border)

# Row is a touch taller to match the new bottom-border row
$ws.Rows.Item(11).RowHeight = 15.75

# --- Selection moves on, like after the user finished entering the row -
$ws.Range("F16").Select() | Out-Null
